{"js": "// Update the worksheet date and every two-digit \u00f7 one-digit division\n// answer cell to the new values from the next day's sheet.\n// Pairs are listed in document order (top-to-bottom, left-to-right) so that\n// the one duplicated \"before\" text (\"89\u00f73=29, 2\") is matched positionally\n// against its own distinct replacement rather than a global find/replace.\nconst replacements = [\n  [\"2024-11-05 Tuesday\", \"2024-11-06 Wednesday\"],\n  [\"34\u00f79=3, 7\", \"81\u00f78=10, 1\"],\n  [\"53\u00f74=13, 1\", \"26\u00f74=6, 2\"],\n  [\"98\u00f75=19, 3\", \"68\u00f76=11, 2\"],\n  [\"89\u00f72=44, 1\", \"14\u00f76=2, 2\"],\n  [\"89\u00f73=29, 2\", \"56\u00f77=8, 0\"],\n  [\"18\u00f75=3, 3\", \"18\u00f79=2, 0\"],\n  [\"76\u00f76=12, 4\", \"22\u00f72=11, 0\"],\n  [\"30\u00f75=6, 0\", \"18\u00f79=2, 0\"],\n  [\"50\u00f79=5, 5\", \"97\u00f77=13, 6\"],\n  [\"71\u00f75=14, 1\", \"92\u00f77=13, 1\"],\n  [\"62\u00f75=12, 2\", \"94\u00f73=31, 1\"],\n  [\"83\u00f74=20, 3\", \"26\u00f74=6, 2\"],\n  [\"36\u00f77=5, 1\", \"54\u00f74=13, 2\"],\n  [\"70\u00f74=17, 2\", \"21\u00f78=2, 5\"],\n  [\"48\u00f72=24, 0\", \"99\u00f77=14, 1\"],\n  [\"77\u00f79=8, 5\", \"69\u00f77=9, 6\"],\n  [\"28\u00f74=7, 0\", \"12\u00f73=4, 0\"],\n  [\"59\u00f76=9, 5\", \"72\u00f75=14, 2\"],\n  [\"29\u00f78=3, 5\", \"66\u00f77=9, 3\"],\n  [\"89\u00f73=29, 2\", \"65\u00f76=10, 5\"],\n  [\"92\u00f72=46, 0\", \"57\u00f77=8, 1\"],\n  [\"73\u00f79=8, 1\", \"15\u00f74=3, 3\"],\n  [\"95\u00f74=23, 3\", \"17\u00f73=5, 2\"],\n  [\"99\u00f75=19, 4\", \"12\u00f72=6, 0\"],\n  [\"65\u00f73=21, 2\", \"65\u00f77=9, 2\"],\n];\n\nconst body = context.document.body;\n\n// Resolve every distinct \"before\" string to its matching ranges FIRST (one\n// search + sync per distinct string, before any text is written), so later\n// writes can never perturb an not-yet-issued search. Occurrences of a\n// duplicated \"before\" string (e.g. \"89\u00f73=29, 2\") come back in document\n// order, matching the top-to-bottom/left-to-right order `replacements`\n// already lists them in.\nconst distinctOld = [...new Set(replacements.map(([oldText]) => oldText))];\nconst foundRanges = new Map();\nfor (const oldText of distinctOld) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  foundRanges.set(oldText, results);\n}\nawait context.sync();\n\n// Now replace each occurrence with its paired replacement text, tracking how\n// many of each distinct \"before\" string have already been consumed.\nconst consumedCount = new Map();\nfor (const [oldText, newText] of replacements) {\n  if (oldText === newText) continue;\n\n  const usedSoFar = consumedCount.get(oldText) ?? 0;\n  const range = foundRanges.get(oldText).items[usedSoFar];\n  consumedCount.set(oldText, usedSoFar + 1);\n\n  range.insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every two-digit \u00f7 one-digit division\n# answer cell to the new values from the next day's sheet.\n# Cells are addressed directly by (row, col) in the single table so the\n# duplicated \"before\" value (\"89\u00f73=29, 2\", which appears twice with two\n# different replacements) is never ambiguous.\n\n$d = $word.ActiveDocument\n\n# Title line above the table: \"2024-11-05 Tuesday\" -> \"2024-11-06 Wednesday\"\n$d.Paragraphs.Item(1).Range.Text = \"2024-11-06 Wednesday\"\n\n$table = $d.Tables.Item(1)\n\n# Word COM Cell(row, col) is 1-based across ALL table rows, including the\n# blank spacer rows between the 5 data rows, so the data rows sit at 1, 5,\n# 9, 13, 17.\n$dataRows = @(1, 5, 9, 13, 17)\n\n$values = @(\n  @(\"81\u00f78=10, 1\", \"26\u00f74=6, 2\", \"68\u00f76=11, 2\", \"14\u00f76=2, 2\", \"56\u00f77=8, 0\"),\n  @(\"18\u00f79=2, 0\", \"22\u00f72=11, 0\", \"18\u00f79=2, 0\", \"97\u00f77=13, 6\", \"92\u00f77=13, 1\"),\n  @(\"94\u00f73=31, 1\", \"26\u00f74=6, 2\", \"54\u00f74=13, 2\", \"21\u00f78=2, 5\", \"99\u00f77=14, 1\"),\n  @(\"69\u00f77=9, 6\", \"12\u00f73=4, 0\", \"72\u00f75=14, 2\", \"66\u00f77=9, 3\", \"65\u00f76=10, 5\"),\n  @(\"57\u00f77=8, 1\", \"15\u00f74=3, 3\", \"17\u00f73=5, 2\", \"12\u00f72=6, 0\", \"65\u00f77=9, 2\")\n)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n  $rowIndex = $dataRows[$r]\n  $rowValues = $values[$r]\n  for ($c = 1; $c -le 5; $c++) {\n    $table.Cell($rowIndex, $c).Range.Text = $rowValues[$c - 1]\n  }\n}\n\nWrite-Output \"done\"\n"}
